$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.181.84"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "'2.014.33"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'225.92"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "'0.608"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'55.27"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "'0.374"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("D10").Value = "'0.0779"
$ws.Range("E10").Value = "  -4.45%  "
$ws.Range("E11").Value = "  -4.05%  "
$ws.Range("D12").Value = "'2.309.98"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "'14.03"
$ws.Range("E13").Value = "  -3.66%  "
$ws.Range("D14").Value = "'19.75"
$ws.Range("E14").Value = "  -3.86%  "
$ws.Range("D15").Value = "'0.736"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "'5.18"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "'2.017.80"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").Value = "'37.068.76"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("D19").Value = "'6.21"
$ws.Range("E19").Value = "  +3.59%  "
$ws.Range("D20").Value = "'68.26"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").Value = "'0.0₃0812"
$ws.Range("E21").Value = "  -4.75%  "
$ws.Range("D22").Value = "'223.93"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").Value = "'2.16"
$ws.Range("E25").Value = "  -4.56%  "
$ws.Range("D26").Value = "'164.33"
$ws.Range("E26").Value = "  -2.25%  "
$ws.Range("D27").Value = "'8.93"
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("D28").Value = "'18.59"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").Value = "'0.124"
$ws.Range("E29").Value = "  -3.92%  "
$ws.Range("D30").Value = "'1.30"
$ws.Range("E30").Value = "  -6.64%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("D32").Value = "'4.40"
$ws.Range("E32").Value = "  -2.44%  "
$ws.Range("D33").Value = "'0.0600"
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D35").Value = "'2.32"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'3.12"
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").Value = "'5.37"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'1.462.54"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("E41").Value = "  -3.44%  "
$ws.Range("D42").Value = "'94.81"
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0911"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.75"
$ws.Range("E44").Value = "  -4.45%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'4.17"
$ws.Range("E45").Value = "  +13.44%  "
$ws.Range("D46").Value = "'15.99"
$ws.Range("E46").Value = "  -4.27%  "
$ws.Range("D47").Value = "'1.12"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'7.05"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D51").Value = "'2.197.24"
$ws.Range("E51").Value = "  -1.44%  "
